# Insert a new "List Bullet" paragraph right after the
# "Docente(s) Responsável(eis) " heading paragraph, containing the two
# teacher names separated by a line break (mirroring the existing
# ListBullet paragraphs elsewhere in this document: first run carries
# the text plus a trailing <w:br/>, second run carries the next line).

$d = $word.ActiveDocument

# Locate the "Docente(s) Responsável(eis)" heading paragraph.
$docenteParagraph = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Docente(s) Responsável(eis)*") {
        $docenteParagraph = $p
        break
    }
}

# Insert a brand-new (empty) paragraph right after it.
$tail = $docenteParagraph.Range
$tail.Collapse(0)            # wdCollapseEnd
$tail.InsertParagraphAfter()

# Grab the freshly created paragraph and give it the "List Bullet" style.
$newParagraph = $docenteParagraph.Next()
$newParagraph.Style = "List Bullet"

# First line of text plus the line break, all in one InsertAfter call so
# they land together inside a single run ending in <w:br/>.
$body = $newParagraph.Range
$body = $d.Range($body.Start, $body.End - 1)   # exclude the paragraph mark
$body.Collapse(0)                              # wdCollapseEnd
$body.InsertAfter("471420 - Carlos Antonio Reis Pereira Baptista" + [char]11)

# Second line goes right after the break, in its own run.
$body = $newParagraph.Range
$pos = $body.End - 1                           # position just before the mark
$tailRange = $d.Range($pos, $pos)
$tailRange.InsertAfter("3586455 - Cassius Olivio Figueiredo Terra Ruchert")
